# Updates the "广州-漫展信息" workbook to the next scrape snapshot:
#  - The "COMICUP 2024SP" (id=83213) listing was removed (cancelled/duplicate entry),
#    so its row is deleted from "展览" and "全部类型", shifting later rows up and
#    renumbering the leading index column.
#  - A few "want to go" counters (column F) ticked up between scrapes for events
#    that are still listed: id=82697 (+1), id=82861 (+2), id=82974 (+3).

function Find-RowByIdSuffix($ws, $idSuffix, $col) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($i = 1; $i -le $lastRow; $i++) {
        $val = $ws.Cells.Item($i, $col).Value()
        if ($val -ne $null -and $val.ToString().EndsWith($idSuffix)) {
            return $i
        }
    }
    return -1
}

function Remove-EventRowAndRenumber($ws, $idSuffix) {
    $row = Find-RowByIdSuffix $ws $idSuffix 8
    if ($row -lt 0) {
        Write-Host "WARNING: row not found for" $idSuffix
        return
    }
    $ws.Rows.Item($row).Delete()
    $lastRow = $ws.UsedRange.Rows.Count
    for ($i = 2; $i -le $lastRow; $i++) {
        $ws.Cells.Item($i, 1).Value = $i - 1
    }
}

function Set-WantToGoCount($ws, $idSuffix, $newValue) {
    $row = Find-RowByIdSuffix $ws $idSuffix 8
    if ($row -lt 0) {
        Write-Host "WARNING: row not found for" $idSuffix
        return
    }
    $ws.Cells.Item($row, 6).Value = $newValue
}

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item(1)    # 展览
$wsShow = $wb.Worksheets.Item(2)    # 演出
$wsLocal = $wb.Worksheets.Item(3)   # 本地生活
$wsAll = $wb.Worksheets.Item(4)     # 全部类型

# Remove the cancelled/duplicate COMICUP 2024SP listing.
Remove-EventRowAndRenumber $wsExpo "id=83213"
Remove-EventRowAndRenumber $wsAll "id=83213"

# Refresh "want to go" counters for listings present in multiple sheets.
Set-WantToGoCount $wsExpo "id=82974" 1214
Set-WantToGoCount $wsAll "id=82974" 1214

Set-WantToGoCount $wsShow "id=82697" 474
Set-WantToGoCount $wsAll "id=82697" 474

Set-WantToGoCount $wsLocal "id=82861" 690
Set-WantToGoCount $wsAll "id=82861" 690
